$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 (5ea61590... file) handoff/handback datetimes refreshed
$wsZhCn.Range("H2").Value = "2016-09-04 06:52:23"
$wsZhCn.Range("K2").Value = "2016-09-04 06:52:41"

# de-de sheet: row 2 (5ea61590... file) handoff/handback datetimes refreshed
$wsDeDe.Range("H2").Value = "2016-09-04 06:52:28"
$wsDeDe.Range("K2").Value = "2016-09-04 06:52:48"

# Overview sheet: column G = "Latest HO Xliff Generate Date"
# Row2 (5ea61590... file) now reflects the new latest de-de handoff datetime
$wsOverview.Range("G2").Value = "2016-09-04 06:52:28"
# Row3 (fa02844e... file) keeps the same displayed value
$wsOverview.Range("G3").Value = "2016-09-04 06:51:20"
